$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (engine quantizes ColumnWidth to 1/6-character grid,
# so 15.666666666666666 is the closest achievable value to the target 16.42578125 width)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666

# Update cell values
$ws.Cells.Item(1, 1).Value = -0.09875120888855804
$ws.Cells.Item(1, 2).Value = 0.09837988052416335
$ws.Cells.Item(2, 1).Value = -0.08006518855109857
$ws.Cells.Item(2, 2).Value = 0.0788492113068795
$ws.Cells.Item(3, 1).Value = -0.0291342931747387
$ws.Cells.Item(3, 2).Value = 0.02881681572984718
$ws.Cells.Item(4, 1).Value = -0.020816815779509668
$ws.Cells.Item(4, 2).Value = 0.02053739502813734
$ws.Cells.Item(5, 1).Value = -0.017537395050394977
$ws.Cells.Item(5, 2).Value = 0.01660129819719014
$ws.Cells.Item(6, 1).Value = -0.006837550604929632
$ws.Cells.Item(6, 2).Value = 0.006483018817716513
$ws.Cells.Item(7, 1).Value = 0.0035169811203581425
$ws.Cells.Item(7, 2).Value = -0.0035994150193703334
$ws.Cells.Item(8, 1).Value = 0.013599414957838007
$ws.Cells.Item(8, 2).Value = -0.013744724093931904
$ws.Cells.Item(9, 1).Value = 0.01574472407898897
$ws.Cells.Item(9, 2).Value = -0.015863564466510827
$ws.Cells.Item(10, 1).Value = 0.017863564453652003
$ws.Cells.Item(10, 2).Value = -0.017870508208996938
$ws.Cells.Item(11, 1).Value = 0.0208705081908942
$ws.Cells.Item(11, 2).Value = -0.02088402653095045
$ws.Cells.Item(12, 1).Value = 0.024384026510830825
$ws.Cells.Item(12, 2).Value = -0.02449896974485366
$ws.Cells.Item(13, 1).Value = -0.017169258079960947
$ws.Cells.Item(13, 2).Value = 0.017081532636124663
$ws.Cells.Item(14, 1).Value = -0.009081532677803317
$ws.Cells.Item(14, 2).Value = 0.00905291216337023
$ws.Cells.Item(15, 1).Value = -0.008052912166294668
$ws.Cells.Item(15, 2).Value = 0.008034435392811012
$ws.Cells.Item(16, 1).Value = -0.006034435401774729
$ws.Cells.Item(16, 2).Value = 0.006003429248790937
$ws.Cells.Item(17, 1).Value = -0.00400342925830266
$ws.Cells.Item(17, 2).Value = 0.003999999979090951
$ws.Cells.Item(18, 1).Value = -0.012314692007223016
$ws.Cells.Item(18, 2).Value = 0.01227189084172764
$ws.Cells.Item(19, 1).Value = -0.008271890864838927
$ws.Cells.Item(19, 2).Value = 0.007978146360796856
$ws.Cells.Item(20, 1).Value = -0.00397814638520444
$ws.Cells.Item(20, 2).Value = 0.0039059698612256
$ws.Cells.Item(21, 1).Value = 0.00009403011419184537
$ws.Cells.Item(21, 2).Value = -0.0001884016048290249
$ws.Cells.Item(22, 1).Value = -0.04571491815864803
$ws.Cells.Item(22, 2).Value = 0.04550078547426928
$ws.Cells.Item(23, 1).Value = -0.04050078550824132
$ws.Cells.Item(23, 2).Value = 0.04009936197727004
$ws.Cells.Item(24, 1).Value = -0.020099362098331852
$ws.Cells.Item(24, 2).Value = 0.019999999877229335
$ws.Cells.Item(25, 1).Value = -0.007763747609724092
$ws.Cells.Item(25, 2).Value = 0.007753500404536595
$ws.Cells.Item(26, 1).Value = -0.005253500424883484
$ws.Cells.Item(26, 2).Value = 0.005242438818090278
$ws.Cells.Item(27, 1).Value = -0.0027424388385717258
$ws.Cells.Item(27, 2).Value = 0.0026876634163399693
$ws.Cells.Item(28, 1).Value = -0.0006876634342054544
$ws.Cells.Item(28, 2).Value = 0.0006619618934493587
$ws.Cells.Item(29, 1).Value = 0.00633803806074873
$ws.Cells.Item(29, 2).Value = -0.006340400914384148
$ws.Cells.Item(30, 1).Value = 0.042255127970753925
$ws.Cells.Item(30, 2).Value = -0.04259652063419095
$ws.Cells.Item(31, 1).Value = 0.045470658721100676
$ws.Cells.Item(31, 2).Value = -0.045593442404589624
$ws.Cells.Item(32, 1).Value = -0.004001117499829121
$ws.Cells.Item(32, 2).Value = 0.003999999976862512
